$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ no = 1;  dir = "C:/TEMP";               file = "controldebug.ini" },
    @{ no = 2;  dir = "C:/TEMP";               file = "friend_list.txt" },
    @{ no = 3;  dir = "C:/TEMP";               file = "SharedDisable.exe" },
    @{ no = 4;  dir = "C:/TEMP";               file = "stop.lnk" },
    @{ no = 5;  dir = "C:/TEMP";               file = "wushowhide.diagcab" },
    @{ no = 6;  dir = "C:/TEMP\TCO_20160823";  file = "data1.cab" },
    @{ no = 7;  dir = "C:/TEMP\TCO_20160823";  file = "data1.hdr" },
    @{ no = 8;  dir = "C:/TEMP\TCO_20160823";  file = "data2.cab" },
    @{ no = 9;  dir = "C:/TEMP\TCO_20160823";  file = "ikernel.ex_" },
    @{ no = 10; dir = "C:/TEMP\TCO_20160823";  file = "layout.bin" },
    @{ no = 11; dir = "C:/TEMP\TCO_20160823";  file = "Setup.exe" },
    @{ no = 12; dir = "C:/TEMP\TCO_20160823";  file = "Setup.ini" },
    @{ no = 13; dir = "C:/TEMP\TCO_20160823";  file = "setup.inx" },
    @{ no = 14; dir = "C:/TEMP\TCO_20160823";  file = "SVRINFO.INI" }
)

# Pass 1: populate column B (directory names) first so the shared-string
# table gets the two directory strings allocated before any filenames,
# matching the order directories/files were first discovered by the scan.
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value2 = $row.dir
    $r = $r + 1
}

# Pass 2: populate the rest (no, d_link formula, file name, f_link formula)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value2 = $row.no
    $ws.Cells.Item($r, 3).Formula = '=HYPERLINK("' + $row.dir + '", "DirView")'
    $ws.Cells.Item($r, 4).Value2 = $row.file
    $ws.Cells.Item($r, 5).Formula = '=HYPERLINK("' + $row.dir + '\' + $row.file + '", "FileView")'
    $r = $r + 1
}
